$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells that are being updated to keep their original text
# representation (with dots as thousands separators / trailing zeros / 
# leading zeros in decimals) instead of being auto-converted to numbers by Excel.
$dCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '26.847.60'
$ws.Range('E2').Value = '  +5.21%  '
$ws.Range('D3').Value = '1.878.17'
$ws.Range('E3').Value = '  +3.95%  '
$ws.Range('D4').Value = '0.9973'
$ws.Range('E4').Value = '  -0.47%  '
$ws.Range('D5').Value = '283.72'
$ws.Range('E5').Value = '  +2.48%  '
$ws.Range('D6').Value = '0.9990'
$ws.Range('E6').Value = '  -0.25%  '
$ws.Range('D7').Value = '0.5192'
$ws.Range('E7').Value = '  +3.50%  '
$ws.Range('D8').Value = '0.3537'
$ws.Range('E8').Value = '  +0.84%  '
$ws.Range('D9').Value = '45.25'
$ws.Range('E9').Value = '  +3.16%  '
$ws.Range('D10').Value = '0.07113'
$ws.Range('E10').Value = '  +6.80%  '
$ws.Range('D11').Value = '20.29'
$ws.Range('E11').Value = '  +1.78%  '
$ws.Range('D12').Value = '0.8233'
$ws.Range('E12').Value = '  -1.59%  '
$ws.Range('D13').Value = '0.07768'
$ws.Range('E13').Value = '  -0.88%  '
$ws.Range('D14').Value = '1.867.09'
$ws.Range('E14').Value = '  +3.33%  '
$ws.Range('D15').Value = '5.180'
$ws.Range('E15').Value = '  +2.55%  '
$ws.Range('D16').Value = '89.95'
$ws.Range('E16').Value = '  +3.19%  '
$ws.Range('D17').Value = '0.9966'
$ws.Range('E17').Value = '  -0.46%  '
$ws.Range('D18').Value = '14.49'
$ws.Range('E18').Value = '  +4.42%  '
$ws.Range('D19').Value = '0.000008173'
$ws.Range('E19').Value = '  +3.42%  '
$ws.Range('D20').Value = '0.9979'
$ws.Range('E20').Value = '  -0.35%  '
$ws.Range('D21').Value = '26.841.27'
$ws.Range('E21').Value = '  +4.83%  '
$ws.Range('E22').Value = '  +1.65%  '
$ws.Range('D23').Value = '10.19'
$ws.Range('E23').Value = '  +2.33%  '
$ws.Range('D24').Value = '6.244'
$ws.Range('E24').Value = '  +3.03%  '
$ws.Range('D25').Value = '2.436'
$ws.Range('E25').Value = '  +15.39%  '
$ws.Range('D26').Value = '145.68'
$ws.Range('E26').Value = '  +3.18%  '
$ws.Range('D27').Value = '17.46'
$ws.Range('E27').Value = '  +3.18%  '
$ws.Range('D28').Value = '1.668'
$ws.Range('E28').Value = '  +0.56%  '
$ws.Range('D29').Value = '111.57'
$ws.Range('E29').Value = '  +2.82%  '
$ws.Range('D30').Value = '4.430'
$ws.Range('E30').Value = '  +3.10%  '
$ws.Range('D31').Value = '4.371'
$ws.Range('E31').Value = '  +3.98%  '
$ws.Range('D32').Value = '0.08867'
$ws.Range('D33').Value = '0.04931'
$ws.Range('E33').Value = '  +2.63%  '
$ws.Range('D34').Value = '1.184'
$ws.Range('E34').Value = '  +5.42%  '
$ws.Range('D35').Value = '0.7495'
$ws.Range('E35').Value = '  +2.31%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = '2.865'
$ws.Range('E36').Value = '  +0.67%  '
$ws.Range('B37').Value = 'MXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D37').Value = '3.291'
$ws.Range('E37').Value = '  +9.02%  '
$ws.Range('D38').Value = '2.431'
$ws.Range('E38').Value = '  +6.07%  '
$ws.Range('D39').Value = '0.5329'
$ws.Range('E39').Value = '  +3.11%  '
$ws.Range('D40').Value = '0.01888'
$ws.Range('E40').Value = '  +1.51%  '
$ws.Range('D41').Value = '0.9763'
$ws.Range('E41').Value = '  +2.12%  '
$ws.Range('D42').Value = '116.59'
$ws.Range('E42').Value = '  +3.60%  '
$ws.Range('D43').Value = '6.319'
$ws.Range('E43').Value = '  +2.47%  '
$ws.Range('D44').Value = '8.217'
$ws.Range('E44').Value = '  +2.45%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = '0.4644'
$ws.Range('E45').Value = '  +1.53%  '
$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').Value = '0.9991'
$ws.Range('E46').Value = '  -0.18%  '
$ws.Range('D47').Value = '0.1374'
$ws.Range('E47').Value = '  -0.43%  '
$ws.Range('D48').Value = '9.513'
$ws.Range('E48').Value = '  +2.84%  '
$ws.Range('D49').Value = '36.81'
$ws.Range('E49').Value = '  +3.18%  '
$ws.Range('D50').Value = '1.522'
$ws.Range('E50').Value = '  +2.20%  '
$ws.Range('E51').Value = '  +1.74%  '
